$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# --- Cells that flip from a numeric value to literal "N/A" style text ---
# (copy format+value from a same-styled source cell in the untouched Murder row,
#  which keeps both the shared-string text and the original style index intact)
$naZero = $ws.Range("C14")   # style 13, text "0"
$naPct  = $ws.Range("E14")   # style 13, text "***.*"

$naZero.Copy($ws.Range("G15"))
$naZero.Copy($ws.Range("C20"))
$naZero.Copy($ws.Range("D20"))
$naZero.Copy($ws.Range("G22"))
$naZero.Copy($ws.Range("G27"))
$naZero.Copy($ws.Range("D29"))
$naZero.Copy($ws.Range("D30"))

$naPct.Copy($ws.Range("H15"))
$naPct.Copy($ws.Range("E20"))
$naPct.Copy($ws.Range("H22"))
$naPct.Copy($ws.Range("H27"))
$naPct.Copy($ws.Range("E29"))
$naPct.Copy($ws.Range("E30"))

# --- Cell that flips from literal "N/A" text back to a normal numeric value ---
$numSrc = $ws.Range("C16")  # style 15 (plain integer count)
$numSrc.Copy($ws.Range("C23"))
$ws.Range("C23").Value = 2

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("L15").Value = -50
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 38
$ws.Range("K16").Value = -28.947368421052
$ws.Range("L16").Value = -18.181818181818
$ws.Range("M16").Value = -25
$ws.Range("N16").Value = -82.580645161290
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -63.636363636363
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 51
$ws.Range("K17").Value = -21.568627450980
$ws.Range("L17").Value = -23.076923076923
$ws.Range("M17").Value = -13.043478260869
$ws.Range("N17").Value = -71.223021582733
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 44
$ws.Range("K18").Value = 29.411764705882
$ws.Range("L18").Value = 7.317073170731
$ws.Range("M18").Value = -12
$ws.Range("N18").Value = -78.217821782178
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -62.5
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -30.434782608695
$ws.Range("I19").Value = 88
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = -11.111111111111
$ws.Range("L19").Value = 20.547945205479
$ws.Range("M19").Value = -10.204081632653
$ws.Range("N19").Value = -5.376344086021
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = -66.666666666666
$ws.Range("M20").Value = -29.032258064516
$ws.Range("N20").Value = -88.717948717948
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = -41.666666666666
$ws.Range("F21").Value = 29
$ws.Range("G21").Value = 52
$ws.Range("H21").Value = -44.230769230769
$ws.Range("I21").Value = 222
$ws.Range("J21").Value = 248
$ws.Range("K21").Value = -10.483870967741
$ws.Range("L21").Value = -5.128205128205
$ws.Range("M21").Value = -15.909090909090
$ws.Range("N21").Value = -72.005044136191
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -55.555555555555
$ws.Range("I23").Value = 42
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = -19.230769230769
$ws.Range("L23").Value = -22.222222222222
$ws.Range("M23").Value = 27.272727272727
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 30
$ws.Range("G24").Value = 59
$ws.Range("H24").Value = -49.152542372881
$ws.Range("I24").Value = 224
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = -10.4
$ws.Range("L24").Value = -24.067796610169
$ws.Range("M24").Value = 8.737864077669
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -83.333333333333
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = -81.25
$ws.Range("I25").Value = 101
$ws.Range("J25").Value = 143
$ws.Range("K25").Value = -29.370629370629
$ws.Range("L25").Value = -38.036809815950
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = -48
$ws.Range("I26").Value = 87
$ws.Range("J26").Value = 78
$ws.Range("K26").Value = 11.538461538461
$ws.Range("L26").Value = 22.535211267605
$ws.Range("M26").Value = -34.586466165413
$ws.Range("L27").Value = -75
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("N29").Value = -86.956521739130
$ws.Range("N30").Value = -88.235294117647
$ws.Range("D31").Value = 2
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = 0

Write-Output "done"